$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's data (daily update) to row 96, mirroring the
# existing date-series pattern in column A (same date number format as
# the cell directly above it, A95).
$ws.Range("A96").Value = 46045
$ws.Range("A96").NumberFormat = $ws.Range("A95").NumberFormat

$ws.Range("B96").Value = 223
$ws.Range("C96").Value = 229
$ws.Range("D96").Value = 219
